{"js": "// Replace the \"three-digit x one-digit\" multiplication prompts in the\n// practice table with the newly generated problems. Each original prompt\n// string is unique within the document, so a simple search + replace per\n// pair is safe and precise.\nconst replacements = [\n  [\"510\u00d77=\", \"636\u00d78=\"],\n  [\"102\u00d78=\", \"107\u00d78=\"],\n  [\"757\u00d78=\", \"430\u00d76=\"],\n  [\"505\u00d74=\", \"370\u00d72=\"],\n  [\"112\u00d79=\", \"296\u00d74=\"],\n  [\"602\u00d77=\", \"500\u00d73=\"],\n  [\"390\u00d72=\", \"368\u00d77=\"],\n  [\"998\u00d73=\", \"484\u00d74=\"],\n  [\"336\u00d78=\", \"385\u00d72=\"],\n  [\"199\u00d75=\", \"374\u00d79=\"],\n  [\"467\u00d76=\", \"475\u00d79=\"],\n  [\"428\u00d73=\", \"476\u00d75=\"],\n  [\"270\u00d73=\", \"193\u00d76=\"],\n  [\"512\u00d74=\", \"564\u00d78=\"],\n  [\"787\u00d72=\", \"248\u00d73=\"],\n  [\"799\u00d75=\", \"738\u00d77=\"],\n  [\"271\u00d78=\", \"315\u00d72=\"],\n  [\"522\u00d74=\", \"114\u00d76=\"],\n  [\"685\u00d74=\", \"147\u00d76=\"],\n  [\"166\u00d74=\", \"393\u00d73=\"],\n  [\"739\u00d77=\", \"822\u00d75=\"],\n  [\"509\u00d78=\", \"807\u00d73=\"],\n  [\"562\u00d78=\", \"691\u00d79=\"],\n  [\"495\u00d75=\", \"611\u00d73=\"],\n  [\"364\u00d76=\", \"222\u00d75=\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Replace the \"three-digit x one-digit\" multiplication prompts in the\n# practice table with the newly generated problems. Each original prompt\n# string is unique within the document, so a Find/Replace pass per pair is\n# safe and precise.\n$d = $word.ActiveDocument\n\n$pairs = @(\n  @(\"510\u00d77=\", \"636\u00d78=\"),\n  @(\"102\u00d78=\", \"107\u00d78=\"),\n  @(\"757\u00d78=\", \"430\u00d76=\"),\n  @(\"505\u00d74=\", \"370\u00d72=\"),\n  @(\"112\u00d79=\", \"296\u00d74=\"),\n  @(\"602\u00d77=\", \"500\u00d73=\"),\n  @(\"390\u00d72=\", \"368\u00d77=\"),\n  @(\"998\u00d73=\", \"484\u00d74=\"),\n  @(\"336\u00d78=\", \"385\u00d72=\"),\n  @(\"199\u00d75=\", \"374\u00d79=\"),\n  @(\"467\u00d76=\", \"475\u00d79=\"),\n  @(\"428\u00d73=\", \"476\u00d75=\"),\n  @(\"270\u00d73=\", \"193\u00d76=\"),\n  @(\"512\u00d74=\", \"564\u00d78=\"),\n  @(\"787\u00d72=\", \"248\u00d73=\"),\n  @(\"799\u00d75=\", \"738\u00d77=\"),\n  @(\"271\u00d78=\", \"315\u00d72=\"),\n  @(\"522\u00d74=\", \"114\u00d76=\"),\n  @(\"685\u00d74=\", \"147\u00d76=\"),\n  @(\"166\u00d74=\", \"393\u00d73=\"),\n  @(\"739\u00d77=\", \"822\u00d75=\"),\n  @(\"509\u00d78=\", \"807\u00d73=\"),\n  @(\"562\u00d78=\", \"691\u00d79=\"),\n  @(\"495\u00d75=\", \"611\u00d73=\"),\n  @(\"364\u00d76=\", \"222\u00d75=\")\n)\n\nforeach ($pair in $pairs) {\n  $find = $d.Content.Find\n  $find.Text = $pair[0]\n  $find.Replacement.Text = $pair[1]\n  $find.Execute($pair[0], $false, $false, $false, $false, $false, $true, 1, $false, $pair[1], 2)\n}\n"}
